$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the date serial value in A1 (45310 -> 45311)
$ws.Range("A1").Value = 45311

# Update "NEGRO" price list (rows 26-32)
$ws.Range("D26").Value = 3450.825
$ws.Range("D27").Value = 4318.805
$ws.Range("D28").Value = 5620.065
$ws.Range("D29").Value = 7941.257
$ws.Range("D30").Value = 11704.373
$ws.Range("D31").Value = 16093.505
$ws.Range("D32").Value = 20689.449

# Update "BLANCO" price list (rows 34-40)
$ws.Range("D34").Value = 3450.825
$ws.Range("D35").Value = 4318.805
$ws.Range("D36").Value = 5620.065
$ws.Range("D37").Value = 7941.257
$ws.Range("D38").Value = 11704.373
$ws.Range("D39").Value = 16093.505
$ws.Range("D40").Value = 20689.449
